$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 17 (pushes old row 17 -> row 18, etc.) ---
$ws.Rows("17").Insert()

# Copy formatting (borders, number formats, fonts) from row 16 into the newly
# inserted blank row 17 so it looks like the other data rows.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the values for the new row 17 (period 2506) ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007126144"
$ws.Range("D17").Value = "VALENTINA DWAN CAMPO PASSO"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 35120
$ws.Range("G17").Value = 878000

# --- Row 16 now represents the newest period, 2507 ---
$ws.Range("E16").Value = "2507"

# --- Row 18 (old row 17, shifted down) now represents period 2505 ---
$ws.Range("E18").Value = "2505"

# --- Update totals: VALOR MORA and Cant. Periodos ---
$ws.Range("E11").Value = 105360
$ws.Range("F13").Value = 3
